# tests/datum_excel_tests.xlsx edit
#
# 1. Switch the active/selected sheet from "Other Tests" to "Gearbox Tests"
#    (clears tabSelected on the old sheet, sets it - plus workbook activeTab -
#    on the new one).
# 2. Move the selection on "Gearbox Tests" from C10 to C5.
# 3. Fill in the (previously empty) measured-value column C for rows 4-10.

$wb  = $excel.ActiveWorkbook
$otherTests   = $wb.Worksheets.Item("Other Tests")
$gearboxTests = $wb.Worksheets.Item("Gearbox Tests")

# Activating Gearbox Tests moves tabSelected/activeTab onto it and off
# Other Tests automatically.
$gearboxTests.Activate()

# New measured values for SURFACE_PAINTED.area, HOUSING.mass, FASTENERS.mass,
# GEARS.mass, DIPSTICK, AIR_NUT and SHAFT_CENTERS respectively.
$gearboxTests.Range("C4").Value  = 3
$gearboxTests.Range("C5").Value  = 2
$gearboxTests.Range("C6").Value  = 4
$gearboxTests.Range("C7").Value  = 5
$gearboxTests.Range("C8").Value  = 3
$gearboxTests.Range("C9").Value  = 6
$gearboxTests.Range("C10").Value = 7

# Leave the cursor on C5, matching the saved selection in the workbook.
$gearboxTests.Range("C5").Select()
